$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-06-22 Sunday"; new = "2025-06-23 Monday"},
    @{old = "96×47="; new = "11×99="},
    @{old = "13×58="; new = "22×11="},
    @{old = "81×82="; new = "30×25="},
    @{old = "95×48="; new = "48×41="},
    @{old = "80×49="; new = "22×33="},
    @{old = "30×65="; new = "70×59="},
    @{old = "86×66="; new = "25×61="},
    @{old = "33×38="; new = "79×78="},
    @{old = "14×41="; new = "95×99="},
    @{old = "16×97="; new = "45×15="},
    @{old = "35×55="; new = "67×27="},
    @{old = "95×83="; new = "23×43="},
    @{old = "71×33="; new = "94×98="},
    @{old = "25×72="; new = "34×37="},
    @{old = "74×77="; new = "93×81="},
    @{old = "70×68="; new = "78×38="},
    @{old = "32×94="; new = "89×35="},
    @{old = "54×23="; new = "25×24="},
    @{old = "43×11="; new = "76×15="},
    @{old = "80×20="; new = "97×55="},
    @{old = "25×45="; new = "43×24="},
    @{old = "12×57="; new = "85×91="},
    @{old = "72×92="; new = "16×96="},
    @{old = "72×50="; new = "63×46="},
    @{old = "96×87="; new = "43×77="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
